# Weekly fruit/vegetable consolidation update:
# insert a new latest-week record at the top (row 2), pushing the
# existing historical rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 2 (first data row),
# shifting all existing data rows down by one.
$ws.Rows.Item(2).Insert()

# The inserted row inherits formatting from the row above (the header);
# clear that so the new data row matches the plain data-row look of the
# rest of the table.
$ws.Rows.Item(2).ClearFormats()

# Restore the date number format on column D (used by every data row).
$ws.Cells.Item(2, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'

# Populate the new row with the latest market record.
$ws.Cells.Item(2, 1).Value = 4
$ws.Cells.Item(2, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(2, 3).Value = 'Los Lagos'
$ws.Cells.Item(2, 4).Value = 44860
$ws.Cells.Item(2, 5).Value = 10
$ws.Cells.Item(2, 6).Value = 300000000
$ws.Cells.Item(2, 7).Value = 'Espárragos'
$ws.Cells.Item(2, 8).Value = 'Sin especificar'
$ws.Cells.Item(2, 9).Value = 'Primera'
$ws.Cells.Item(2, 10).Value = 200
$ws.Cells.Item(2, 11).Value = 1700
$ws.Cells.Item(2, 12).Value = 1700
$ws.Cells.Item(2, 13).Value = 1700
$ws.Cells.Item(2, 14).Value = '$/kilo'
$ws.Cells.Item(2, 15).Value = 'Provincia de Linares'
$ws.Cells.Item(2, 16).Value = 1700
$ws.Cells.Item(2, 17).Value = 1
$ws.Cells.Item(2, 18).Value = 'Hortaliza'
